$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling / wording of the "Meaning" column values.
$ws.Range("E3").Value = "Instruction Set Architecture"
$ws.Range("E7").Value = "Arithmetic  Logic Unit"
$ws.Range("E8").Value = "Static Random Access  Memory"
$ws.Range("E10").Value = "Serial Peripheral  Interface"

# Update the view state: scroll so row 13 is the top row and select E10.
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollRow = 13

$wb.Save()
